$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Formula = "=B2*100"
$ws.Range("F2").Formula = "=C2"

$ws.Range("E3:E6").Formula = "=B3*100"
$ws.Range("F3:F6").Formula = "=C3"

$ws.Range("E2:E6").NumberFormat = "0.000"
$ws.Range("F2:F6").NumberFormat = "0.0000"

$ws.Range("F6").Select()

$ws.PageSetup.Orientation = 1
